# Apply the 2024-05-02 GitHub Actions cryptos-list refresh: updated
# prices / 1h-volume percentages, plus the PEPE/Mantle row-34<->35 swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.378.04'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '2.998.88'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '''563.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '''138.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.65%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").Value = '2.984.66'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").Value = '''0.0000230'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").Value = '''33.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("E16").Value = '  +7.62%  '
$ws.Range("D17").Value = '3.492.99'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '2.991.97'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '59.297.48'
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").Value = '''430.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").Value = '''13.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.80%  '
$ws.Range("D22").Value = '''0.720'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.20%  '
$ws.Range("D23").Value = '''7.15'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = '''13.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = '''81.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +7.97%  '
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("D30").Value = '''7.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = '''25.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("D33").Value = '''0.0988'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.97%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '''0.990'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.62%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0768'
$ws.Range("E35").Value = '  +11.64%  '
$ws.Range("E36").Value = '  +3.15%  '
$ws.Range("E37").Value = '  -3.67%  '
$ws.Range("D38").Value = '''49.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").Value = '''2.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("D41").Value = '''402.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.38%  '
$ws.Range("D42").Value = '2.771.28'
$ws.Range("E42").Value = '  +3.67%  '
$ws.Range("D43").Value = '''0.0352'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("E44").Value = '  -2.21%  '
$ws.Range("D45").Value = '''0.251'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.76%  '
$ws.Range("D47").Value = '''34.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +17.50%  '
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '''120.68'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '''2.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").Value = '''23.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.86%  '
